$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 0.004732036689895636
$ws.Range("H2").Value = 0.01031643294823183
$ws.Range("K2").Value = 5.046707586018868
$ws.Range("L2").Value = "[1.3410462271884764, 8.752368944849259]"
$ws.Range("M2").Value = 0.007954304748214902
$ws.Range("N2").Value = 0.007954304748214902
$ws.Range("O2").Value = -1.371105502467618
$ws.Range("P2").Value = "[-2.3145267197618495, -0.42768428517338597]"
$ws.Range("Q2").Value = 0.004693095350625009
$ws.Range("R2").Value = 0.004693095350625009
$ws.Range("S2").Value = 10.31570430230792
$ws.Range("T2").Value = "[8.178006853028524, 12.453401751587318]"
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 4.800800800800801
$ws.Range("X2").Value = 1.497497497497498
$ws.Range("Y2").Value = 8.104104104104103

# Row 3 updates
$ws.Range("E3").Value = 23.22000000000019
$ws.Range("G3").Value = [double]"4.029251813308932e-07"
$ws.Range("H3").Value = [double]"5.968530519634601e-06"
$ws.Range("K3").Value = 6.137480332485705
$ws.Range("L3").Value = "[3.2547617308354866, 9.020198934135923]"
$ws.Range("M3").Value = [double]"3.71401079894973e-05"
$ws.Range("N3").Value = [double]"7.428021597899459e-05"
$ws.Range("O3").Value = -1.056631763369539
$ws.Range("P3").Value = "[-1.54721079636254, -0.5660527303765388]"
$ws.Range("Q3").Value = [double]"3.030341751264487e-05"
$ws.Range("R3").Value = [double]"6.060683502528974e-05"
$ws.Range("S3").Value = 10.24227839968147
$ws.Range("T3").Value = "[8.692301469126829, 11.79225533023612]"
$ws.Range("W3").Value = 3.904864864864894
$ws.Range("X3").Value = 2.091891891891906
$ws.Range("Y3").Value = 5.717837837837882
